$wb = $excel.ActiveWorkbook

# --- Original sheet ("Sheet") ---
$shSheet = $wb.Worksheets.Item("Sheet")

# --- Add the 5 new worksheets, in order, after the last existing sheet ---
$shAndy = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$shAndy.Name = "Andy"

$shZY = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$shZY.Name = "ZY"

$shDash = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$shDash.Name = "Dash"

$shPC = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$shPC.Name = "PC"

$shWH = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$shWH.Name = "WH"

# --- Populate each new sheet with the same header row (A1:E1) as "Sheet",
#     preserving values + styles via a direct range copy ---
$shSheet.Range("A1:E1").Copy($shAndy.Range("A1:E1"))
$shSheet.Range("A1:E1").Copy($shZY.Range("A1:E1"))
$shSheet.Range("A1:E1").Copy($shDash.Range("A1:E1"))
$shSheet.Range("A1:E1").Copy($shPC.Range("A1:E1"))
$shSheet.Range("A1:E1").Copy($shWH.Range("A1:E1"))

# --- Column widths on "Andy" ---
$shAndy.Columns.Item(2).ColumnWidth = 37.666666666666664
$shAndy.Columns.Item(3).ColumnWidth = 29.833333333333332
$shAndy.Columns.Item(4).ColumnWidth = 30.666666666666668
$shAndy.Columns.Item(5).ColumnWidth = 31.833333333333332

# --- Per-sheet selections / views ---

# Andy: active cell B2
$shAndy.Activate()
$shAndy.Range("B2").Select()

# ZY: whole column B selected, active cell B1
$shZY.Activate()
$shZY.Columns.Item(2).Select()

# Dash: A1:E1 selected
$shDash.Activate()
$shDash.Range("A1:E1").Select()

# PC: A1:E1 selected
$shPC.Activate()
$shPC.Range("A1:E1").Select()

# "Sheet": scroll so column D is the leftmost visible column, select A1:E1
$shSheet.Activate()
$shSheet.Range("A1:E1").Select()
$excel.ActiveWindow.ScrollColumn = 4

# WH: A1:E1 selected; leave as the active/selected tab (matches activeTab=5,
# tabSelected on the last sheet)
$shWH.Activate()
$shWH.Range("A1:E1").Select()
